# Weekly update: insert two new daily-price rows (most recent date 44438)
# for "Ciboulette" at the top of the dated block (row 140), pushing the
# existing rows 140:221 down to 142:223.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 140; everything below
# (old rows 140-221) shifts down to 142-223, carrying its formatting
# (incl. the date style on column D) along with it.
$ws.Rows("140:141").Insert()

# --- New row 140 ---
$ws.Cells.Item(140, 1).Value = 9
$ws.Cells.Item(140, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(140, 3).Value = "Metropolitana"
$ws.Cells.Item(140, 4).Value = 44438
$ws.Cells.Item(140, 5).Value = 13
$ws.Cells.Item(140, 6).Value = 100112039
$ws.Cells.Item(140, 7).Value = "Ciboulette"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 250
$ws.Cells.Item(140, 11).Value = 3000
$ws.Cells.Item(140, 12).Value = 3300
$ws.Cells.Item(140, 13).Value = 3150
$ws.Cells.Item(140, 14).Value = "`$/docena de atados"
$ws.Cells.Item(140, 15).Value = "Región Metropolitana"
$ws.Cells.Item(140, 16).Value = 1050
$ws.Cells.Item(140, 17).Value = 3
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# --- New row 141 ---
$ws.Cells.Item(141, 1).Value = 9
$ws.Cells.Item(141, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(141, 3).Value = "Metropolitana"
$ws.Cells.Item(141, 4).Value = 44438
$ws.Cells.Item(141, 5).Value = 13
$ws.Cells.Item(141, 6).Value = 100112039
$ws.Cells.Item(141, 7).Value = "Ciboulette"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Segunda"
$ws.Cells.Item(141, 10).Value = 160
$ws.Cells.Item(141, 11).Value = 2500
$ws.Cells.Item(141, 12).Value = 2800
$ws.Cells.Item(141, 13).Value = 2650
$ws.Cells.Item(141, 14).Value = "`$/docena de atados"
$ws.Cells.Item(141, 15).Value = "Región Metropolitana"
$ws.Cells.Item(141, 16).Value = 883
$ws.Cells.Item(141, 17).Value = 3
$ws.Cells.Item(141, 18).Value = "Hortaliza"
